$d = $word.ActiveDocument

function Insert-ParaXml($range, $innerXml) {
    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes" ?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/part.xml" pkg:contentType="application/xml">
<pkg:xmlData>
$innerXml
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $range.InsertXML($xml)
}

$rPrTNR = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr>'
$pPrLvl1 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="480" w:lineRule="auto"/><w:textAlignment w:val="baseline"/>' + $rPrTNR + '</w:pPr>'

# --- 1. DreamHost paragraph: split " for hosting using .xyz" into
#        " for hosting using ." + "xyz", and wrap DreamHost / xyz with
#        spellcheck proofErr markers.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "DreamHost for hosting using .xyz*") {
        $target = $p
        break
    }
}
$innerXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
$pPrLvl1
<w:proofErr w:type="spellStart"/>
<w:r>$rPrTNR<w:t>DreamHost</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>$rPrTNR<w:t xml:space="preserve"> for hosting using .</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r>$rPrTNR<w:t>xyz</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
"@
Insert-ParaXml $target.Range $innerXml

# --- 2. "WP User Registration for user accounts" -> "Ultimate Member plug-in for user accounts"
#        (split across two runs)
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "WP User Registration for user accounts*") {
        $target = $p
        break
    }
}
$innerXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
$pPrLvl1
<w:r>$rPrTNR<w:t>Ultimate Member plug-in</w:t></w:r>
<w:r>$rPrTNR<w:t xml:space="preserve"> for user accounts</w:t></w:r>
</w:p>
"@
Insert-ParaXml $target.Range $innerXml

# --- 3. "Amazon Product Advertising API to display products from Amazon"
#        becomes TWO paragraphs:
#        a) "Dropshipping & Affiliation with Amazon to display products from Amazon"
#           (replacing the original paragraph, with proofErr around "Dropshipping")
#        b) a brand-new paragraph inserted right after:
#           "Amazon Product Advertising API to keep products from Amazon updated"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Amazon Product Advertising API to display products from Amazon*") {
        $target = $p
        break
    }
}
$innerXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
$pPrLvl1
<w:proofErr w:type="spellStart"/>
<w:r>$rPrTNR<w:t>Dropshipping</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>$rPrTNR<w:t xml:space="preserve"> &amp; Affiliation with Amazon to display products from Amazon</w:t></w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
$pPrLvl1
<w:r>$rPrTNR<w:t xml:space="preserve">Amazon Product Advertising API to </w:t></w:r>
<w:r>$rPrTNR<w:t>keep</w:t></w:r>
<w:r>$rPrTNR<w:t xml:space="preserve"> products from Amazon</w:t></w:r>
<w:r>$rPrTNR<w:t xml:space="preserve"> updated</w:t></w:r>
</w:p>
"@
Insert-ParaXml $target.Range $innerXml
